$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja1: replace the December OT rows with the new March OT rows, then drop
# the now-empty trailing row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A2").Value = 43915
$ws1.Range("B2").Value = 43915
$ws1.Range("C2").Value = "MPI_CUATRIMESTRAL_03-20_PUERTO BERRIO"
$ws1.Range("D2").Value = "GD042"

$ws1.Range("A3").Value = 43901
$ws1.Range("B3").Value = 43901
$ws1.Range("C3").Value = "MPI_SEMESTRAL_03-20_SAN JOSE DEL NUS"
$ws1.Range("D3").Value = "GD234"

$ws1.Range("A4").Value = 43903
$ws1.Range("B4").Value = 43903
$ws1.Range("C4").Value = "MPI_SEMESTRAL_03-20_LA DANTA"
$ws1.Range("D4").Value = "GD352"

$ws1.Range("A5:D5").ClearContents()

$ws1.Rows.Item(81).Delete()

$ws1.Columns.Item(4).ColumnWidth = 26.6328125
$ws1.Range("C23").Select()

# The AutoFilter memory left on Hoja1 now points at a much smaller range.
$fd = $wb.Names.Item("Hoja1!_FilterDatabase")
$fd.RefersTo = '=Hoja1!$A$1:$D$14'

# ---------------------------------------------------------------------------
# Hoja3: the whole reference table gets wiped (values only - formatting and
# merged cells stay put).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Range("A1:G5").ClearContents()
$ws3.Range("C10").Select()

# ---------------------------------------------------------------------------
# Hoja2: same three OT rows (site/ID/start/end/task), row 5 emptied, and the
# sheet grown down to row 40 of blank bordered rows. Column D becomes visible.
# Handled last so Hoja2 ends up the active sheet/tab, as in the source file.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")

$ws2.Range("A2").Value = "S2230"
$ws2.Range("B2").Value = "GD042"
$ws2.Range("C2").Value = 5297292
$ws2.Range("D2").Value = 5297294
$ws2.Range("E2").Value = "MPI_CUATRIMESTRAL_03-20_PUERTO BERRIO"

$ws2.Range("A3").Value = "S8404"
$ws2.Range("B3").Value = "GD234"
$ws2.Range("C3").Value = 5297298
$ws2.Range("D3").Value = 5297300
$ws2.Range("E3").Value = "MPI_SEMESTRAL_03-20_SAN JOSE DEL NUS"

$ws2.Range("A4").Value = "S9796"
$ws2.Range("B4").Value = "GD352"
$ws2.Range("C4").Value = 5297301
$ws2.Range("D4").Value = 5297302
$ws2.Range("E4").Value = "MPI_SEMESTRAL_03-20_LA DANTA"

$ws2.Range("A5:E5").ClearContents()

$ws2.Range("A6:E40").Borders.LineStyle = 1
$ws2.Columns.Item(4).Hidden = $false
$ws2.Columns.Item(4).ColumnWidth = 11.54296875
$ws2.Range("C4").Select()
